$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Days Commit in GitHub" row (row 9) ---
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = ""
$ws.Rows.Item(9).AutoFit()

# --- Update "Numbers of Commits in GitHub" row (row 10) ---
$ws.Range("D10").Value = 40
$ws.Range("E10").Value = ""
$ws.Rows.Item(10).AutoFit()

# --- Update "Delete Ad" row (row 43) ---
$ws.Range("C43").Value = "Yes"
$ws.Range("E43").Value = "Service problem"

# --- Update "Implement Paging" comment (row 35) ---
$ws.Range("E35").Value = "paging is not work properly"

# --- Update "Change Password" row (row 47) ---
$ws.Range("C47").Value = "Yes"
$ws.Range("E47").Value = "Service problem"

# --- Update "Admin List Towns" row (row 51) ---
$ws.Range("C51").Value = "Yes"
$ws.Range("E51").Value = "Service problem"

# --- Update sheet view: selection moves to D12 (also clears stale topLeftCell scroll state) ---
$ws.Range("D12").Select()
